$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 65 (Excel copies formatting down from the row above,
# so the new row inherits row 64's style automatically).
$ws.Rows(65).Insert()

# New GLOF record: Rembesdalskåka / Demmevatnet, 2025 event.
$ws.Range("A65").Value = "Rembesdalskåka"
$ws.Range("B65").Value = "Demmevatnet"
$ws.Range("C65").Value = 45774
$ws.Range("D65").Value = "28.12.2024-27.04.2025"
$ws.Range("E65").Value = "no"
$ws.Range("F65").Value = "60.5425599"
$ws.Range("G65").Value = "7.3190236"
$ws.Range("H65").Value = "1239"
$ws.Range("I65").Value = 2025
$ws.Range("J65").Value = "April"

# The trailing sorted block (previously A109:J112, sorted by date) shifted
# down by one row because of the inserted row above it; refresh the sort
# bookkeeping to match (data is already in order, so this just re-records
# the range/condition without moving anything).
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("C110:C113")) | Out-Null
$sort.SetRange($ws.Range("A110:J113"))
$sort.Header = 2
$sort.Apply()

# Match the recorded view state after the edit.
$ws.Range("D65").Select()
$excel.ActiveWindow.ScrollRow = 31
